$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update installed/demand load calculations for the correct engine power
$ws.Range("D6").Value = 2.76
$ws.Range("E6").Value = 2

$ws.Range("D7").Value = 1.42

$ws.Range("D8").Value = 6.65

$ws.Range("D9").Value = 17.73

$ws.Range("D10").Value = 0.9

# Update active selection on the sheet to reflect the reviewed data range
$ws.Range("A2:F11").Select()
